# Scheduled-runner refresh of market-board derived figures (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets. Values below are plain numbers (no formulas in this workbook), so we
# just poke the refreshed figures straight into the affected cells per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1180.6
$ws.Range("I5").Value = 803.6667
$ws.Range("K5").Value = 803.6667
$ws.Range("M5").Value = -688.6667
$ws.Range("H17").Value = 1091.3448
$ws.Range("J17").Value = 1075.4036
$ws.Range("L17").Value = 3226.2108
$ws.Range("N17").Value = -3562.2108
$ws.Range("H28").Value = 487.04166
$ws.Range("I28").Value = 350
$ws.Range("K28").Value = 350
$ws.Range("M28").Value = 135
$ws.Range("H113").Value = 4380.3335
$ws.Range("I113").Value = 4298
$ws.Range("K113").Value = 4298
$ws.Range("M113").Value = -1044
$ws.Range("H116").Value = 6432.067
$ws.Range("J116").Value = 4534.4287
$ws.Range("L116").Value = 4534.4287
$ws.Range("N116").Value = -11418.4287
$ws.Range("H130").Value = 104716.4
$ws.Range("I130").Value = 65000
$ws.Range("J130").Value = 114645.5
$ws.Range("K130").Value = 65000
$ws.Range("L130").Value = 114645.5
$ws.Range("M130").Value = -59980
$ws.Range("N130").Value = -124685.5
$ws.Range("H132").Value = 7879.0425
$ws.Range("I132").Value = 4959.355
$ws.Range("K132").Value = 14878.065
$ws.Range("M132").Value = -12348.065
$ws.Range("H140").Value = 98162
$ws.Range("I140").Value = 47349
$ws.Range("K140").Value = 47349
$ws.Range("M140").Value = -42169

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1022568.25
$ws.Range("J2").Value = 2423.818
$ws.Range("L2").Value = 2423.818
$ws.Range("N2").Value = -2649.818
$ws.Range("H32").Value = 3334.519
$ws.Range("I32").Value = 2937.352
$ws.Range("K32").Value = 2937.352
$ws.Range("M32").Value = -2650.352
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H116").Value = 1022568.25
$ws.Range("J116").Value = 2423.818
$ws.Range("L116").Value = 2423.818
$ws.Range("N116").Value = -7011.818

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1022568.25
$ws.Range("J3").Value = 2423.818
$ws.Range("L3").Value = 2423.818
$ws.Range("N3").Value = -2651.818
$ws.Range("H134").Value = 3426.889
$ws.Range("I134").Value = 3426.889
$ws.Range("K134").Value = 10280.667
$ws.Range("M134").Value = -7745.667000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 187.5
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -88
$ws.Range("N4").Value = -374
$ws.Range("H22").Value = 1862.3334
$ws.Range("J22").Value = 1730.2858
$ws.Range("L22").Value = 1730.2858
$ws.Range("N22").Value = -2430.2858
$ws.Range("H99").Value = 6574.551
$ws.Range("J99").Value = 6229.5493
$ws.Range("L99").Value = 6229.5493
$ws.Range("N99").Value = -9225.549299999999
$ws.Range("H122").Value = 3629.6099
$ws.Range("I122").Value = 3454.4
$ws.Range("K122").Value = 10363.2
$ws.Range("M122").Value = -7913.200000000001
$ws.Range("H126").Value = 6574.551
$ws.Range("J126").Value = 6229.5493
$ws.Range("L126").Value = 18688.6479
$ws.Range("N126").Value = -23628.6479
$ws.Range("H138").Value = 113431.71
$ws.Range("J138").Value = 113431.71
$ws.Range("L138").Value = 113431.71
$ws.Range("N138").Value = -123711.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 646.1539
$ws.Range("I5").Value = 468.95456
$ws.Range("K5").Value = 1406.86368
$ws.Range("M5").Value = -1294.86368
$ws.Range("H92").Value = 637.1818
$ws.Range("I92").Value = 672
$ws.Range("K92").Value = 2016
$ws.Range("M92").Value = -768
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").ClearContents()
$ws.Range("N101").Value = 0
$ws.Range("H114").Value = 878.75
$ws.Range("I114").Value = 1020.8571
$ws.Range("J114").Value = 768.2222
$ws.Range("K114").Value = 3062.5713
$ws.Range("L114").Value = 2304.6666
$ws.Range("M114").Value = 191.4287000000004
$ws.Range("N114").Value = -8812.6666
$ws.Range("H135").Value = 646.1539
$ws.Range("I135").Value = 468.95456
$ws.Range("K135").Value = 4220.59104
$ws.Range("M135").Value = -1685.59104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20003646
$ws.Range("I80").Value = 27062382
$ws.Range("J80").Value = 3895.6667
$ws.Range("K80").Value = 27062382
$ws.Range("L80").Value = 3895.6667
$ws.Range("M80").Value = -27061384
$ws.Range("N80").Value = -5891.6667
$ws.Range("H83").Value = 20003646
$ws.Range("I83").Value = 27062382
$ws.Range("J83").Value = 3895.6667
$ws.Range("K83").Value = 135311910
$ws.Range("L83").Value = 19478.3335
$ws.Range("M83").Value = -135306918
$ws.Range("N83").Value = -29462.3335
$ws.Range("H107").Value = 427.5
$ws.Range("I107").Value = 397.22223
$ws.Range("K107").Value = 397.22223
$ws.Range("M107").Value = 1522.77777
$ws.Range("H132").Value = 4773.961
$ws.Range("I132").Value = 4494.7666
$ws.Range("K132").Value = 13484.2998
$ws.Range("M132").Value = -10954.2998
$ws.Range("H140").Value = 56435.57
$ws.Range("J140").Value = 63341.5
$ws.Range("L140").Value = 63341.5
$ws.Range("N140").Value = -73701.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1967.6
$ws.Range("I7").Value = 1885.875
$ws.Range("K7").Value = 1885.875
$ws.Range("M7").Value = -1773.875
$ws.Range("H46").Value = 3464.5
$ws.Range("J46").Value = 4749.1
$ws.Range("L46").Value = 4749.1
$ws.Range("N46").Value = -5125.1
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H126").Value = 1967.6
$ws.Range("I126").Value = 1885.875
$ws.Range("K126").Value = 5657.625
$ws.Range("M126").Value = -3187.625
$ws.Range("H132").Value = 17002.69
$ws.Range("I132").Value = 19826.7
$ws.Range("K132").Value = 59480.10000000001
$ws.Range("M132").Value = -56950.10000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 12999.5
$ws.Range("I2").Value = 6999
$ws.Range("J2").Value = 19000
$ws.Range("K2").Value = 6999
$ws.Range("L2").Value = 19000
$ws.Range("M2").Value = -6887
$ws.Range("N2").Value = -19224
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("N34").Value = 0
$ws.Range("H37").Value = 12800
$ws.Range("I37").Value = 9000
$ws.Range("K37").Value = 9000
$ws.Range("M37").Value = -8797
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1851
$ws.Range("H42").Value = 61000
$ws.Range("I42").Value = 45000
$ws.Range("K42").Value = 45000
$ws.Range("M42").Value = -44622
$ws.Range("H43").Value = 35499.75
$ws.Range("I43").Value = 30999.666
$ws.Range("K43").Value = 30999.666
$ws.Range("M43").Value = -30850.666
$ws.Range("H107").Value = 2893
$ws.Range("I107").Value = 2674.4
$ws.Range("J107").Value = 3986
$ws.Range("K107").Value = 8023.200000000001
$ws.Range("L107").Value = 11958
$ws.Range("M107").Value = -6103.200000000001
$ws.Range("N107").Value = -15798
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0
